$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.508.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.96%  '

$ws.Range("D3").Value = "'2.313.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.01%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = "'310.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.72%  '

$ws.Range("D6").Value = "'105.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.92%  '

$ws.Range("E7").Value = '  +1.49%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = "'0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.35%  '

$ws.Range("D10").Value = "'37.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.98%  '

$ws.Range("D11").Value = "'52.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.54%  '

$ws.Range("D12").Value = "'0.0816"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.79%  '

$ws.Range("E13").Value = '  -0.80%  '

$ws.Range("D14").Value = "'7.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.99%  '

$ws.Range("D15").Value = "'2.672.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.96%  '

$ws.Range("E16").Value = '  +3.96%  '

$ws.Range("D17").Value = "'2.307.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.02%  '

$ws.Range("E18").Value = '  +3.57%  '

$ws.Range("D19").Value = "'43.413.28"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = "'12.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").Value = "'0.0₃0932"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.13%  '

$ws.Range("E22").Value = '  +3.99%  '

$ws.Range("D23").Value = "'68.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.12%  '

$ws.Range("D24").Value = "'242.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.70%  '

$ws.Range("E25").Value = '  +3.03%  '

$ws.Range("D26").Value = "'2.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.20%  '

$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("E28").Value = '  +5.50%  '

$ws.Range("E29").Value = '  +12.14%  '

$ws.Range("D30").Value = "'37.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.82%  '

$ws.Range("D31").Value = "'9.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.28%  '

$ws.Range("D32").Value = "'166.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.95%  '

$ws.Range("D33").Value = "'5.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.48%  '

$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D34").Value = "'18.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.12%  '

$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("E36").Value = '  +6.70%  '

$ws.Range("E37").Value = '  +1.99%  '

$ws.Range("D38").Value = "'3.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.65%  '

$ws.Range("D39").Value = "'4.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.85%  '

$ws.Range("E40").Value = '  +4.18%  '

$ws.Range("E41").Value = '  +2.88%  '

$ws.Range("E42").Value = '  +0.71%  '

$ws.Range("E43").Value = '  +20.24%  '

$ws.Range("E44").Value = '  +4.04%  '

$ws.Range("D45").Value = "'1.998.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.58%  '

$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = "'3.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.76%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'19.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.98%  '

$ws.Range("D48").Value = "'10.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.81%  '

$ws.Range("D49").Value = "'57.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.00%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = "'1.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.32%  '

$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").Value = "'2.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.72%  '
